$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data region (rows 1-18, cols A-N) so we can rewrite it
# cleanly with the updated / reordered content (an extra "trend_epi" rds row
# was inserted, and rows 11-18 were reshuffled).
$ws.Range("A1:N18").ClearContents()

# Full target grid: row, column, value
$cellData = @(
      @(1,1,"nm"),
      @(1,2,"sql_str"),
      @(1,3,"adjust_func"),
      @(1,4,"dir"),
      @(1,5,"suffix"),
      @(1,6,"fn"),
      @(1,7,"saving_func"),
      @(1,8,"Sun"),
      @(1,9,"Mon"),
      @(1,10,"Tue"),
      @(1,11,"Wed"),
      @(1,12,"Thu"),
      @(1,13,"Fri"),
      @(1,14,"Sat"),
      @(2,1,"web_epi"),
      @(2,2,"select * from all_cases_web;"),
      @(2,3,"keep_only_web_epi_cols"),
      @(2,4,"//Ncr-a_irbv2s/irbv2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DATABASE"),
      @(2,5," "),
      @(2,6,"all_cases_web_current.xlsx"),
      @(2,13,"X"),
      @(3,1,"web_epi"),
      @(3,2,"select * from all_cases_web;"),
      @(3,3,"keep_only_web_epi_cols"),
      @(3,4,"//Ncr-a_irbv2s/irbv2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DATABASE"),
      @(3,5," "),
      @(3,6,"all_cases_web_current.rds"),
      @(3,7,"saveRDS"),
      @(3,13,"X"),
      @(4,1,"trend_epi"),
      @(4,2,"select * from all_cases;"),
      @(4,3,"keep_only_trend_epi_cols"),
      @(4,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/EPI SUMMARY/Trend analysis/_Current/_Source Data/CaseReportForm"),
      @(4,5," "),
      @(4,6,"basic_CRF_extract.xlsx"),
      @(4,9,"X"),
      @(4,12,"X"),
      @(5,1,"trend_epi"),
      @(5,2,"select * from all_cases;"),
      @(5,3,"keep_only_trend_epi_cols"),
      @(5,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/EPI SUMMARY/Trend analysis/_Current/_Source Data/CaseReportForm"),
      @(5,5," "),
      @(5,6,"basic_CRF_extract.rds"),
      @(5,7,"saveRDS"),
      @(5,9,"X"),
      @(5,12,"X"),
      @(6,1,"Dashboard"),
      @(6,2,"select * from all_cases;"),
      @(6,4,"//Ncr-a_irbv2s/irbv2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DASHBOARD"),
      @(6,5,"format(Sys.Date(), '%Y-%m-%d')"),
      @(6,6,"qry_allcases_current.xlsx"),
      @(6,8,"X"),
      @(6,9,"X"),
      @(6,10,"X"),
      @(6,11,"X"),
      @(6,12,"X"),
      @(6,13,"X"),
      @(6,14,"X"),
      @(7,1,"Dashboard"),
      @(7,2,"select * from all_cases;"),
      @(7,4,"//Ncr-a_irbv2s/irbv2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DASHBOARD"),
      @(7,5,"format(Sys.Date(), '%Y-%m-%d')"),
      @(7,6,"qry_allcases_current.rds"),
      @(7,7,"saveRDS"),
      @(7,8,"X"),
      @(7,9,"X"),
      @(7,10,"X"),
      @(7,11,"X"),
      @(7,12,"X"),
      @(7,13,"X"),
      @(7,14,"X"),
      @(8,1,"epi"),
      @(8,2,"select * from all_cases;"),
      @(8,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/SAS_Analysis/Domestic data"),
      @(8,5,"format(Sys.Date(), '%Y-%m-%d')"),
      @(8,6,"qry_allcases {suffix}_DISCOVER.xlsx"),
      @(8,8,"X"),
      @(8,9,"X"),
      @(8,10,"X"),
      @(8,11,"X"),
      @(8,12,"X"),
      @(8,13,"X"),
      @(8,14,"X"),
      @(9,1,"epi"),
      @(9,2,"select * from all_cases;"),
      @(9,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/SAS_Analysis/Domestic data"),
      @(9,5,"format(Sys.Date(), '%Y-%m-%d')"),
      @(9,6,"qry_allcases {suffix}_DISCOVER.rds"),
      @(9,7,"saveRDS"),
      @(9,8,"X"),
      @(9,9,"X"),
      @(9,10,"X"),
      @(9,11,"X"),
      @(9,12,"X"),
      @(9,13,"X"),
      @(9,14,"X"),
      @(10,1,"HCDaily"),
      @(10,2,"select * from data_hub;"),
      @(10,3,"remove_pt_cols"),
      @(10,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/Data Requests/SituationalAwareness Dashboard Daily Extract"),
      @(10,5,"format(Sys.Date() ,'%Y%m%d')"),
      @(10,6,"{suffix}_HCDaily_DISCOVER.xlsx"),
      @(10,8,"X"),
      @(10,9,"X"),
      @(10,10,"X"),
      @(10,11,"X"),
      @(10,12,"X"),
      @(10,13,"X"),
      @(10,14,"X"),
      @(11,1,"HCDaily"),
      @(11,2,"select * from data_hub;"),
      @(11,3,"remove_pt_cols"),
      @(11,4,"L:/HPOC/Active Events/001-20 COVID-19/Dashboard"),
      @(11,5,"format(Sys.Date() ,'%Y%m%d')"),
      @(11,6,"{suffix}_HCDaily_DISCOVER.xlsx"),
      @(11,8,"X"),
      @(11,9,"X"),
      @(11,10,"X"),
      @(11,11,"X"),
      @(11,12,"X"),
      @(11,13,"X"),
      @(11,14,"X"),
      @(12,1,"STATCAN"),
      @(12,2,"select * from statscan where classification='confirmed';"),
      @(12,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/STATCAN"),
      @(12,5,"format(Sys.Date() ,'%Y%m%d')"),
      @(12,6,"Weekly Extended Dataset_{suffix}_DISCOVER.xlsx"),
      @(12,8,"X"),
      @(13,1,"WHO"),
      @(13,2,"select * from who;"),
      @(13,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DATABASE/OUTPUTS/WHO"),
      @(13,5,"format(Sys.Date() ,'%d%B%Y')"),
      @(13,6,"Canada_COVID19_WHO_linelist-{suffix}_DISCOVER.xlsx"),
      @(13,12,"X"),
      @(14,1,"WHO"),
      @(14,2,"select * from who;"),
      @(14,4,"L:/HPOC/Active Events/001-20 COVID-19/Operations/Surveillance - Epi. Diagnostics/Canada_COVID19_WHO_linelist"),
      @(14,5,"format(Sys.Date() ,'%d%B%Y')"),
      @(14,6,"Canada_COVID19_WHO_linelist-{suffix}_DISCOVER.xlsx"),
      @(14,12,"X"),
      @(15,1,"modeling"),
      @(15,2,"select * from modelling_data;"),
      @(15,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/PHAC MODELLING/Domestic data"),
      @(15,5,"format(Sys.Date() ,'%Y-%m-%d')"),
      @(15,6,"Domestic surveillance data - {suffix}_DISCOVER.xlsx"),
      @(15,8,"X"),
      @(15,9,"X"),
      @(15,10,"X"),
      @(15,11,"X"),
      @(15,12,"X"),
      @(15,13,"X"),
      @(15,14,"X"),
      @(16,1,"modeling"),
      @(16,2,"select * from modelling_data;"),
      @(16,4,"//Ncr-a-phacc1s/phacc1/HPOC/Active Events/001-20 COVID-19/Operations/Emerging Science/Modelling Group/Models and Data/Domestic Surveillance Data"),
      @(16,5,"format(Sys.Date() ,'%Y-%m-%d')"),
      @(16,6,"Domestic surveillance data - {suffix}_DISCOVER.xlsx"),
      @(16,8,"X"),
      @(16,9,"X"),
      @(16,10,"X"),
      @(16,11,"X"),
      @(16,12,"X"),
      @(16,13,"X"),
      @(16,14,"X"),
      @(17,1,"modeling"),
      @(17,2,"select * from modelling_data;"),
      @(17,4,"https://storphacidpcbns02.blob.core.windows.net/hcdaily"),
      @(17,5,"format(Sys.Date() ,'%Y-%m-%d')"),
      @(17,6,"Domestic surveillance data - {suffix}_DISCOVER.xlsx"),
      @(17,7,"save_azure"),
      @(17,8,"X"),
      @(17,9,"X"),
      @(17,10,"X"),
      @(17,11,"X"),
      @(17,12,"X"),
      @(17,13,"X"),
      @(17,14,"X"),
      @(18,1,"datahub"),
      @(18,2,"select * from all_cases;"),
      @(18,3,"make_data_hub"),
      @(18,4,"https://storhpocnspdatalakeprod.blob.core.windows.net/hcdaily/data"),
      @(18,5," "),
      @(18,6,"current_DataHub_DISCOVER.xlsx"),
      @(18,7,"save_azure"),
      @(18,8,"X"),
      @(18,9,"X"),
      @(18,10,"X"),
      @(18,11,"X"),
      @(18,12,"X"),
      @(18,13,"X"),
      @(18,14,"X"),
      @(19,1,"countSummary"),
      @(19,2,"select * from all_cases;"),
      @(19,3,"count_summary"),
      @(19,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DATABASE/OUTPUTS/CountSummary"),
      @(19,5,"format(Sys.Date() ,'%Y-%m-%d')"),
      @(19,6,"CountSummary_{suffix}.xlsx"),
      @(19,8,"X"),
      @(19,9,"X"),
      @(19,10,"X"),
      @(19,11,"X"),
      @(19,12,"X"),
      @(19,13,"X"),
      @(19,14,"X"),
      @(20,1,"db_errs"),
      @(20,2,"select * from all_cases;"),
      @(20,3,"db_error_report_by_case"),
      @(20,4,"//Ncr-a_irbv2s/IRBV2/PHAC/IDPCB/CIRID/VIPS-SAR/EMERGENCY PREPAREDNESS AND RESPONSE HC4/EMERGENCY EVENT/WUHAN UNKNOWN PNEU - 2020/DATA AND ANALYSIS/DATABASE/DISCOVER/Data Quality/db_errs"),
      @(20,5,"format(Sys.Date() ,'%Y-%m-%d')"),
      @(20,6,"db_errs {suffix}.xlsx"),
      @(20,8,"X"),
      @(20,9,"X"),
      @(20,10,"X"),
      @(20,11,"X"),
      @(20,12,"X"),
      @(20,13,"X"),
      @(20,14,"X")

)

foreach ($cell in $cellData) {
    $r = $cell[0]
    $c = $cell[1]
    $v = $cell[2]
    $ws.Cells.Item($r, $c).Value = $v
}

# Update the selected cell shown when the workbook is opened.
$ws.Range("A17").Select()

# The new "modeling" row (row 17) carries a hyperlink on its directory path.
$ws.Hyperlinks.Add($ws.Range("D17"), "https://storphacidpcbns02.blob.core.windows.net/hcdaily")
